# Auto-generated PowerShell-style Excel COM-interop edit script
# Applies numeric corrections to the F/G "interest count" and "min price" columns
# across the four worksheets (展览, 演出, 本地生活, 全部类型), matching the commit diff.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet2 = $wb.Worksheets.Item(2)   # 演出
$sheet3 = $wb.Worksheets.Item(3)   # 本地生活
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 ---
$sheet1.Range("F2").Value = 1852
$sheet1.Range("G4").Value = "不可售"
$sheet1.Range("F5").Value = 52
$sheet1.Range("G6").Value = 78
$sheet1.Range("F8").Value = 191
$sheet1.Range("F9").Value = 627
$sheet1.Range("F11").Value = 467
$sheet1.Range("F12").Value = 748
$sheet1.Range("F13").Value = 1453
$sheet1.Range("F14").Value = 1213
$sheet1.Range("F15").Value = 1443
$sheet1.Range("F16").Value = 25
$sheet1.Range("F17").Value = 1236
$sheet1.Range("F18").Value = 299
$sheet1.Range("F19").Value = 1586
$sheet1.Range("F20").Value = 775
$sheet1.Range("F22").Value = 323
$sheet1.Range("F25").Value = 1369
$sheet1.Range("F27").Value = 74
$sheet1.Range("F28").Value = 799
$sheet1.Range("F30").Value = 1061
$sheet1.Range("F32").Value = 971
$sheet1.Range("F35").Value = 1319
$sheet1.Range("F36").Value = 33
$sheet1.Range("F38").Value = 1063
$sheet1.Range("F39").Value = 16
$sheet1.Range("F40").Value = 43
$sheet1.Range("F41").Value = 36
$sheet1.Range("F43").Value = 1589
$sheet1.Range("F44").Value = 92
$sheet1.Range("F46").Value = 790
# --- 演出 ---
$sheet2.Range("F4").Value = 131
$sheet2.Range("F11").Value = 1427
$sheet2.Range("F14").Value = 2527
$sheet2.Range("F15").Value = 1184
$sheet2.Range("F18").Value = 218
$sheet2.Range("F23").Value = 433
$sheet2.Range("F27").Value = 0
$sheet2.Range("F34").Value = 131
$sheet2.Range("F46").Value = 124
$sheet2.Range("F47").Value = 54
# --- 本地生活 ---
$sheet3.Range("F5").Value = 2754
$sheet3.Range("F6").Value = 4518
$sheet3.Range("F7").Value = 121
$sheet3.Range("F9").Value = 535
$sheet3.Range("F10").Value = 645
$sheet3.Range("F12").Value = 207
$sheet3.Range("F13").Value = 783
$sheet3.Range("F14").Value = 195
$sheet3.Range("F15").Value = 441
# --- 全部类型 ---
$sheet4.Range("F2").Value = 1852
$sheet4.Range("F4").Value = 2754
$sheet4.Range("F6").Value = 4518
$sheet4.Range("F7").Value = 645
$sheet4.Range("F8").Value = 52
$sheet4.Range("F9").Value = 207
$sheet4.Range("F10").Value = 207
$sheet4.Range("F11").Value = 783
$sheet4.Range("F12").Value = 195
$sheet4.Range("G13").Value = 78
$sheet4.Range("F15").Value = 191
$sheet4.Range("F16").Value = 1427
$sheet4.Range("F17").Value = 627
$sheet4.Range("F18").Value = 467
$sheet4.Range("F19").Value = 748
$sheet4.Range("F20").Value = 2527
$sheet4.Range("F21").Value = 1184
$sheet4.Range("F22").Value = 1453
$sheet4.Range("F23").Value = 1213
$sheet4.Range("F24").Value = 1443
$sheet4.Range("F25").Value = 1236
$sheet4.Range("F26").Value = 218
$sheet4.Range("F27").Value = 299
$sheet4.Range("F29").Value = 1586
$sheet4.Range("F30").Value = 775
$sheet4.Range("F31").Value = 323
$sheet4.Range("F32").Value = 441
$sheet4.Range("F33").Value = 433
$sheet4.Range("F34").Value = 1369
$sheet4.Range("F36").Value = 799
$sheet4.Range("F38").Value = 1061
$sheet4.Range("F40").Value = 971
$sheet4.Range("F42").Value = 1063
$sheet4.Range("F44").Value = 1589
$sheet4.Range("F45").Value = 92
$sheet4.Range("F47").Value = 790
$sheet4.Range("F51").Value = 54

Write-Output "Applied all cell updates."
